# Daily attendance processing - 2025-10-30 07:43:00
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Widen the "Status"-adjacent "Recorded By" helper column (col I, index 9): 10 -> 14 ---
$ws.Columns.Item(9).ColumnWidth = 13.14

# --- 2) Reorder "Recorded By" (column G) name lists: move "System" right after any leading
#        lowercase "system" token (or to the front otherwise) ---
$gUpdates = @{
    2 = 'system, System, backup@backdoor.com'
    3 = 'System, dnasr281@gmail.com'
    4 = 'System, backup@backdoor.com'
    5 = 'System, backup@backdoor.com'
    6 = 'System, dnasr281@gmail.com'
    7 = 'System, admin@admin.com'
    8 = 'System, backup@backdoor.com'
    10 = 'System, dnasr281@gmail.com'
    11 = 'System, dnasr281@gmail.com'
    12 = 'System, dnasr281@gmail.com'
    13 = 'System, dnasr281@gmail.com'
    14 = 'System, dnasr281@gmail.com'
    15 = 'System, dnasr281@gmail.com'
    17 = 'System, dnasr281@gmail.com'
    18 = 'System, dnasr281@gmail.com'
    19 = 'System, dnasr281@gmail.com'
    20 = 'System, dnasr281@gmail.com'
    21 = 'System, dnasr281@gmail.com'
    22 = 'System, dnasr281@gmail.com'
    24 = 'System, dnasr281@gmail.com'
    29 = 'system, System, backup@backdoor.com'
    30 = 'System, dnasr281@gmail.com'
    31 = 'System, backup@backdoor.com'
    32 = 'System, backup@backdoor.com'
    33 = 'System, dnasr281@gmail.com'
    34 = 'System, admin@admin.com'
    35 = 'System, backup@backdoor.com'
    37 = 'System, dnasr281@gmail.com'
    38 = 'System, dnasr281@gmail.com'
    39 = 'System, dnasr281@gmail.com'
    40 = 'System, dnasr281@gmail.com'
    41 = 'System, dnasr281@gmail.com'
    42 = 'System, dnasr281@gmail.com'
    44 = 'System, dnasr281@gmail.com'
    45 = 'System, dnasr281@gmail.com'
    46 = 'System, dnasr281@gmail.com'
    47 = 'System, dnasr281@gmail.com'
    48 = 'System, dnasr281@gmail.com'
    49 = 'System, dnasr281@gmail.com'
    51 = 'System, dnasr281@gmail.com'
    56 = 'system, System, backup@backdoor.com'
    57 = 'System, dnasr281@gmail.com'
    58 = 'System, backup@backdoor.com'
    59 = 'System, backup@backdoor.com'
    60 = 'System, dnasr281@gmail.com'
    61 = 'System, admin@admin.com'
    62 = 'System, backup@backdoor.com'
    64 = 'System, dnasr281@gmail.com'
    65 = 'System, dnasr281@gmail.com'
    66 = 'System, dnasr281@gmail.com'
    67 = 'System, dnasr281@gmail.com'
    68 = 'System, dnasr281@gmail.com'
    69 = 'System, dnasr281@gmail.com'
    71 = 'System, dnasr281@gmail.com'
    72 = 'System, dnasr281@gmail.com'
    73 = 'System, dnasr281@gmail.com'
    74 = 'System, dnasr281@gmail.com'
    75 = 'System, dnasr281@gmail.com'
    76 = 'System, dnasr281@gmail.com'
    78 = 'System, dnasr281@gmail.com'
    83 = 'System, backup@backdoor.com'
    84 = 'System, backup@backdoor.com'
    85 = 'System, backup@backdoor.com'
    86 = 'System, dnasr281@gmail.com'
    87 = 'System, dnasr281@gmail.com'
    88 = 'System, dnasr281@gmail.com'
    89 = 'System, dnasr281@gmail.com'
    93 = 'System, dnasr281@gmail.com'
    95 = 'System, dnasr281@gmail.com'
    96 = 'System, dnasr281@gmail.com'
    97 = 'System, dnasr281@gmail.com'
    99 = 'System, dnasr281@gmail.com'
    102 = 'System, dnasr281@gmail.com'
    104 = 'System, dnasr281@gmail.com'
    109 = 'System, backup@backdoor.com'
    110 = 'System, backup@backdoor.com'
    111 = 'System, backup@backdoor.com'
    112 = 'System, dnasr281@gmail.com'
    113 = 'System, dnasr281@gmail.com'
    114 = 'System, dnasr281@gmail.com'
    115 = 'System, dnasr281@gmail.com'
    119 = 'System, dnasr281@gmail.com'
    121 = 'System, dnasr281@gmail.com'
    122 = 'System, dnasr281@gmail.com'
    123 = 'System, dnasr281@gmail.com'
    125 = 'System, dnasr281@gmail.com'
    128 = 'System, dnasr281@gmail.com'
    130 = 'System, dnasr281@gmail.com'
    135 = 'System, backup@backdoor.com'
    136 = 'System, backup@backdoor.com'
    137 = 'System, backup@backdoor.com'
    138 = 'System, dnasr281@gmail.com'
    139 = 'System, dnasr281@gmail.com'
    140 = 'System, dnasr281@gmail.com'
    141 = 'System, dnasr281@gmail.com'
    145 = 'System, dnasr281@gmail.com'
    147 = 'System, dnasr281@gmail.com'
    148 = 'System, dnasr281@gmail.com'
    149 = 'System, dnasr281@gmail.com'
    151 = 'System, dnasr281@gmail.com'
    154 = 'System, dnasr281@gmail.com'
    156 = 'System, dnasr281@gmail.com'
}
foreach ($row in $gUpdates.Keys) {
    $ws.Range("G" + $row).Value = $gUpdates[$row]
}

# --- 3) Mark the three still-not-recorded 30/10/2025 sessions (B2D/B2E/B2F, session 23) as
#        "Not Recorded" (was "Pending") and recolor them red/pink like the other "Not Recorded" rows ---
foreach ($rownum in 105, 131, 157) {
    $range = $ws.Range("A" + $rownum + ":I" + $rownum)
    $range.Interior.Color = 12695295
    $range.Font.Color = 0
    $ws.Range("I" + $rownum).Value = "Not Recorded"
}

# --- 4) Update the summary statistics affected by the 3 sessions moving from Pending to Missing ---
$ws.Range("L7").Value = 3    # Missing Sessions: 0 -> 3
$ws.Range("L8").Value = 18   # Pending Sessions: 21 -> 18

# Group Statistics rows for B2D (18), B2E (19) and B2F (20): Missing 0 -> 1, Pending 4 -> 3
foreach ($rownum in 18, 19, 20) {
    $ws.Range("P" + $rownum).Value = 1
    $ws.Range("Q" + $rownum).Value = 3
}

